$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2048929663608563
$ws.Range("C2").Value = 0.5351681957186545
$ws.Range("J2").Value = 0.009174311926605505
$ws.Range("P2").Value = 0.1498470948012232
$ws.Range("S2").Value = 0.1009174311926606
$ws.Range("B3").Value = 0.01104972375690608
$ws.Range("C3").Value = 0.03867403314917127
$ws.Range("J3").Value = 0.03867403314917127
$ws.Range("P3").Value = 0.7016574585635359
$ws.Range("S3").Value = 0.2099447513812155
$ws.Range("J4").Value = 0.04651162790697674
$ws.Range("P4").Value = 0.6511627906976745
$ws.Range("S4").Value = 0.3023255813953488
$ws.Range("B6").Value = 0.06343283582089553
$ws.Range("D6").Value = 0.007462686567164179
$ws.Range("E6").Value = 0.003731343283582089
$ws.Range("F6").Value = 0.05970149253731343
$ws.Range("J6").Value = 0.376865671641791
$ws.Range("O6").Value = 0.01119402985074627
$ws.Range("Q6").Value = 0.1529850746268657
$ws.Range("R6").Value = 0.04104477611940299
$ws.Range("S6").Value = 0.2835820895522388
$ws.Range("B7").Value = 0.1226415094339623
$ws.Range("D7").Value = 0.03773584905660377
$ws.Range("F7").Value = 0.07075471698113207
$ws.Range("J7").Value = 0.1273584905660377
$ws.Range("O7").Value = 0.009433962264150943
$ws.Range("Q7").Value = 0.2169811320754717
$ws.Range("R7").Value = 0.07075471698113207
$ws.Range("S7").Value = 0.3443396226415094
$ws.Range("B8").Value = 0.09981515711645102
$ws.Range("D8").Value = 0.011090573012939
$ws.Range("E8").Value = 0.001848428835489834
$ws.Range("F8").Value = 0.05730129390018484
$ws.Range("J8").Value = 0.11090573012939
$ws.Range("O8").Value = 0.01478743068391867
$ws.Range("Q8").Value = 0.1959334565619224
$ws.Range("R8").Value = 0.09057301293900184
$ws.Range("S8").Value = 0.4177449168207024
$ws.Range("B9").Value = 0.1038251366120219
$ws.Range("D9").Value = 0.03278688524590164
$ws.Range("F9").Value = 0.04371584699453552
$ws.Range("J9").Value = 0.1038251366120219
$ws.Range("O9").Value = 0.03278688524590164
$ws.Range("Q9").Value = 0.2513661202185792
$ws.Range("R9").Value = 0.09836065573770492
$ws.Range("S9").Value = 0.3333333333333333
$ws.Range("B10").Value = 0.1109337589784517
$ws.Range("D10").Value = 0.01755786113328013
$ws.Range("E10").Value = 0.0007980845969672786
$ws.Range("F10").Value = 0.0710295291300878
$ws.Range("J10").Value = 0.1109337589784517
$ws.Range("O10").Value = 0.009577015163607342
$ws.Range("Q10").Value = 0.2130885873902634
$ws.Range("R10").Value = 0.0742218675179569
$ws.Range("S10").Value = 0.3918595371109337
$ws.Range("G11").Value = 0.1438127090301003
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.1906354515050167
$ws.Range("L11").Value = 0.5785953177257525
$ws.Range("S11").Value = 0.01003344481605351
$ws.Range("G12").Value = 0.7853107344632768
$ws.Range("J12").Value = 0.1694915254237288
$ws.Range("L12").Value = 0.01694915254237288
$ws.Range("S12").Value = 0.02824858757062147
$ws.Range("G13").Value = 0.5806451612903226
$ws.Range("J13").Value = 0.3870967741935484
$ws.Range("S13").Value = 0.03225806451612903
$ws.Range("F15").Value = 0.02487562189054726
$ws.Range("H15").Value = 0.1741293532338309
$ws.Range("I15").Value = 0.07960199004975124
$ws.Range("J15").Value = 0.373134328358209
$ws.Range("K15").Value = 0.03980099502487562
$ws.Range("M15").Value = 0.01492537313432836
$ws.Range("O15").Value = 0.05970149253731343
$ws.Range("S15").Value = 0.2338308457711443
$ws.Range("F16").Value = 0.02463054187192118
$ws.Range("H16").Value = 0.2660098522167488
$ws.Range("I16").Value = 0.06403940886699508
$ws.Range("J16").Value = 0.3694581280788177
$ws.Range("K16").Value = 0.07881773399014778
$ws.Range("M16").Value = 0.02463054187192118
$ws.Range("O16").Value = 0.07389162561576355
$ws.Range("S16").Value = 0.09852216748768473
$ws.Range("F17").Value = 0.03607214428857716
$ws.Range("H17").Value = 0.2064128256513026
$ws.Range("I17").Value = 0.09619238476953908
$ws.Range("J17").Value = 0.3947895791583166
$ws.Range("K17").Value = 0.1122244488977956
$ws.Range("M17").Value = 0.03006012024048096
$ws.Range("O17").Value = 0.04208416833667335
$ws.Range("S17").Value = 0.08216432865731463
$ws.Range("F18").Value = 0.03783783783783784
$ws.Range("H18").Value = 0.1891891891891892
$ws.Range("I18").Value = 0.0918918918918919
$ws.Range("J18").Value = 0.3945945945945946
$ws.Range("K18").Value = 0.1189189189189189
$ws.Range("M18").Value = 0.02702702702702703
$ws.Range("O18").Value = 0.05405405405405406
$ws.Range("S18").Value = 0.08648648648648649
$ws.Range("F19").Value = 0.03855421686746988
$ws.Range("H19").Value = 0.2530120481927711
$ws.Range("I19").Value = 0.07309236947791165
$ws.Range("J19").Value = 0.3301204819277108
$ws.Range("K19").Value = 0.1124497991967871
$ws.Range("M19").Value = 0.02891566265060241
$ws.Range("O19").Value = 0.07068273092369477
$ws.Range("S19").Value = 0.09317269076305221

Write-Host "Applied 108 cell updates"
